$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update training hours in column F ---
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = 2
$ws.Range("F10").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 1

$ws.Range("F15").Value = 2
$ws.Range("F16").Value = 2
$ws.Range("F17").Value = 2
$ws.Range("F18").Value = 2
$ws.Range("F19").Value = 2
$ws.Range("F20").Value = 2
$ws.Range("F21").Value = 2
$ws.Range("F22").Value = 2
$ws.Range("F23").Value = 2
$ws.Range("F24").Value = 2

# --- Align formatting of F15:F24 with F7:F14 (right-aligned / wrap, thin border) ---
$ws.Range("F7").Copy()
$ws.Range("F15:F24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update Team Skills Audit row 24 (PowerPoint/presentation -> Documentation) ---
$ws.Range("B24").Value = "Documentation"
$ws.Range("E24").Value = "To make professional documentation for our website."

# --- Update view state (scroll position / selection) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F24").Select()
